$d = $word.ActiveDocument

# 1) Header cell "pvalues" -> "p"
$d.Content.Find.Execute("pvalues", $true, $false, $false, $false, $false,
                         $true, 1, $false, "p", 2)

# 2) p-values in the table body
$d.Content.Find.Execute(".205", $true, $false, $false, $false, $false,
                         $true, 1, $false, ".386", 2)
$d.Content.Find.Execute(".119", $true, $false, $false, $false, $false,
                         $true, 1, $false, ".386", 2)
$d.Content.Find.Execute(".232", $true, $false, $false, $false, $false,
                         $true, 1, $false, ".386", 2)
$d.Content.Find.Execute(".480", $true, $false, $false, $false, $false,
                         $true, 1, $false, ".599", 2)

# 3) Add a new paragraph after "Dependent Variable: var5"
$p = $d.Paragraphs.Last
$p.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "Multiple tests correction applied to p values: Benjamini-Hochberg"
